$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-20 from 45184 to 45185
$ws.Range("C2:C20").Value = 45185
